# Updates cryptos list price (D) and volume-change (E) columns to latest values.
# Generated from the authoritative diff of xl/worksheets/sheet1.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are plain numeric-looking strings (e.g. "227.60",
# "0.0693"); force text format first so Excel keeps them as exact strings
# instead of normalizing them into floating point numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.330.70'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.801.86'
$ws.Range('E3').Value = '  +0.78%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.60'
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.578'
$ws.Range('E6').Value = '  +3.97%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '36.05'
$ws.Range('E8').Value = '  +9.37%  '
$ws.Range('E9').Value = '  +2.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0693'
$ws.Range('E10').Value = '  +0.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0965'
$ws.Range('E11').Value = '  +2.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.061.97'
$ws.Range('E12').Value = '  +0.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.55'
$ws.Range('E13').Value = '  +3.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.775.22'
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.645'
$ws.Range('E15').Value = '  +1.56%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.50'
$ws.Range('E16').Value = '  +4.81%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '34.319.18'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.06'
$ws.Range('E18').Value = '  +0.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.54'
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0795'
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.50'
$ws.Range('E21').Value = '  +2.36%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.18'
$ws.Range('E23').Value = '  +0.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '172.18'
$ws.Range('E24').Value = '  +2.10%  '
$ws.Range('E25').Value = '  +3.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.90'
$ws.Range('E26').Value = '  +7.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.86'
$ws.Range('E27').Value = '  +1.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.119'
$ws.Range('E28').Value = '  +2.90%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.05'
$ws.Range('E30').Value = '  +0.82%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0532'
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.85'
$ws.Range('E32').Value = '  +1.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.25'
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('E34').Value = '  +0.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.394.21'
$ws.Range('E35').Value = '  -1.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.674'
$ws.Range('E36').Value = '  -1.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.47'
$ws.Range('E37').Value = '  -5.11%  '
$ws.Range('E38').Value = '  -0.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0190'
$ws.Range('E39').Value = '  -0.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.23'
$ws.Range('E40').Value = '  +11.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.960'
$ws.Range('E41').Value = '  +2.40%  '
$ws.Range('E42').Value = '  +1.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '82.04'
$ws.Range('E43').Value = '  -2.87%  '
$ws.Range('E44').Value = '  +0.44%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.58'
$ws.Range('E45').Value = '  -3.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.03'
$ws.Range('E46').Value = '  -0.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0501'
$ws.Range('E47').Value = '  -5.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.962.56'
$ws.Range('E48').Value = '  +0.85%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '104.81'
$ws.Range('E49').Value = '  -0.59%  '
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('E51').Value = '  -0.23%  '
